$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on sheet Hoja1 (A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 4.55 = 17968.18 pesos`n✅ 17968.18 pesos = 4.51 = 934.76 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update rate figures on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 220
$wsTasas.Range("O10").Value = 3953
$wsTasas.Range("N12").Value = 3980
$wsTasas.Range("O12").Value = 207.051
